$wb = $excel.ActiveWorkbook

# Add the new "Comments" header to the relevant sheets and update their
# selection/view state to match the authored edit.

$ws2 = $wb.Worksheets.Item("Withdraw History")
$ws2.Range("E1").Value = "Comments"

$ws3 = $wb.Worksheets.Item("Deposit History")
$ws3.Range("E1").Value = "Comments"

$ws4 = $wb.Worksheets.Item("Transfer History")
$ws4.Range("E1").Value = "Comments"

$ws5 = $wb.Worksheets.Item("Absolute History")
$ws5.Range("E1").Value = "Comments"

# Restore selections on the non-active sheets first (selecting a range also
# activates its sheet, so the last selection made below determines the
# workbook's final active tab).
$ws5.Range("E5").Select()
$ws3.Range("E1").Select()
$ws4.Range("E1").Select()

# "Withdraw History" ends up as the active sheet/tab.
$ws2.Range("E1").Select()
